$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from the row above for the styled columns (A and E),
# then overwrite with the new values.
$ws.Range("A45").Copy($ws.Range("A46"))
$ws.Range("E45").Copy($ws.Range("E46"))

$ws.Range("A46").Value = 45
$ws.Range("B46").Value = "gibraltar"
$ws.Range("C46").Value = "national-league"
$ws.Range("D46").Value = "2023-2024"
$ws.Range("E46").Value = 45281.875
$ws.Range("F46").Value = "Lincoln Red Imps"
$ws.Range("G46").Value = 1
$ws.Range("H46").Value = "St Josephs"
$ws.Range("I46").Value = 2
$ws.Range("J46").Value = 1.83
$ws.Range("K46").Value = "21/12/2023 10:16"
$ws.Range("L46").Value = 1.68
$ws.Range("M46").Value = "21/12/2023 20:02"
$ws.Range("N46").Value = 3.91
$ws.Range("O46").Value = "21/12/2023 10:16"
$ws.Range("P46").Value = 3.81
$ws.Range("Q46").Value = "21/12/2023 20:51"
$ws.Range("R46").Value = 3.22
$ws.Range("S46").Value = "21/12/2023 10:16"
$ws.Range("T46").Value = 4
$ws.Range("U46").Value = "21/12/2023 20:51"
$ws.Range("V46").Value = "https://www.betexplorer.com/football/gibraltar/national-league/lincoln-red-imps-st-josephs/pxvMGxSN/"
